# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to Text format before writing, so that
# values such as "1.002" or "0.9985" are preserved exactly as strings
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.071.43"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.816.44"
$ws.Range("E3").Value = "  +2.19%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "337.62"
$ws.Range("E5").Value = "  -0.57%  "

# Row 6 - USDC
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.37%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4267"
$ws.Range("E7").Value = "  +11.63%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3507"
$ws.Range("E8").Value = "  +2.66%  "

# Row 9 - OKB
$ws.Range("D9").Value = "45.62"
$ws.Range("E9").Value = "  -2.56%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "1.149"
$ws.Range("E10").Value = "  +0.65%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.07449"
$ws.Range("E11").Value = "  +0.83%  "

# Row 12 - Solana
$ws.Range("D12").Value = "23.02"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13 - BinanceUSD
$ws.Range("D13").Value = "1.0000"
$ws.Range("E13").Value = "  -0.18%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.269"
$ws.Range("E14").Value = "  -1.69%  "

# Row 15 - was Chainlink, now WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.817.17"
$ws.Range("E15").Value = "  +2.24%  "

# Row 16 - was WrappedEther, now Chainlink
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.284"
$ws.Range("E16").Value = "  -1.86%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001085"
$ws.Range("E17").Value = "  +0.95%  "

# Row 18 - TRON
$ws.Range("D18").Value = "0.06681"
$ws.Range("E18").Value = "  +0.28%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "82.02"
$ws.Range("E19").Value = "  -0.52%  "

# Row 20 - Dai (D unchanged)
$ws.Range("E20").Value = "  +0.05%  "

# Row 21 - was Uniswap, now Avalanche
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "17.29"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22 - was Avalanche, now Uniswap
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.463"
$ws.Range("E22").Value = "  +0.96%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "28.099.75"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "12.00"
$ws.Range("E24").Value = "  -0.57%  "

# Row 25 - Toncoin (D unchanged)
$ws.Range("E25").Value = "  +0.42%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.495"
$ws.Range("E26").Value = "  +3.53%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "20.75"
$ws.Range("E27").Value = "  +0.17%  "

# Row 28 - Monero (D unchanged)
$ws.Range("E28").Value = "  +1.30%  "

# Row 29 - WrappedliquidstakedEther2.0 (E unchanged)
$ws.Range("D29").Value = "2.018.78"

# Row 30 - ImmutableX
$ws.Range("D30").Value = "1.304"
$ws.Range("E30").Value = "  -9.85%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "132.62"
$ws.Range("E31").Value = "  -1.24%  "

# Row 32 - HuobiToken
$ws.Range("D32").Value = "4.057"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.967"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "0.09196"
$ws.Range("E34").Value = "  +3.28%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "12.38"
$ws.Range("E35").Value = "  -2.61%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.02362"
$ws.Range("E36").Value = "  -1.54%  "

# Row 37 - TheSandbox
$ws.Range("D37").Value = "0.6741"
$ws.Range("E37").Value = "  -1.45%  "

# Row 38 - InternetComputer(DFINITY)
$ws.Range("D38").Value = "5.253"
$ws.Range("E38").Value = "  -0.73%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "0.06276"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40 - Algorand
$ws.Range("D40").Value = "0.2169"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - WEMIXTOKEN
$ws.Range("D41").Value = "1.494"
$ws.Range("E41").Value = "  -0.35%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.220"
$ws.Range("E42").Value = "  -1.50%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "8.175"
$ws.Range("E43").Value = "  -0.57%  "

# Row 44 - Frax
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "14.09"
$ws.Range("E45").Value = "  -1.37%  "

# Row 46 - PancakeSwap (D unchanged)
$ws.Range("E46").Value = "  +0.17%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.6134"
$ws.Range("E47").Value = "  -2.09%  "

# Row 48 - Quant
$ws.Range("D48").Value = "128.77"
$ws.Range("E48").Value = "  -3.20%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "2.047"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50 - EOS
$ws.Range("D50").Value = "1.178"
$ws.Range("E50").Value = "  -2.53%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.07107"
$ws.Range("E51").Value = "  -5.24%  "
